# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Nectarín - Super Queen, Especial/Primera)
# right after the existing row 362, shifting the rest of the table down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 363 and 364 (pushes old 363.. down to 365..)
$ws.Range("A363:A364").EntireRow.Insert()

# --- Row 363: Nectarín / Super Queen / Especial ---
$ws.Cells.Item(363, 1).Value  = 4
$ws.Cells.Item(363, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(363, 3).Value  = "Los Lagos"
$ws.Cells.Item(363, 4).Value  = 44585
$ws.Cells.Item(363, 5).Value  = 10
$ws.Cells.Item(363, 6).Value  = "Fruta"
$ws.Cells.Item(363, 7).Value  = 100103
$ws.Cells.Item(363, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(363, 9).Value  = 100103006
$ws.Cells.Item(363, 10).Value = "Nectarín"
$ws.Cells.Item(363, 11).Value = "Super Queen"
$ws.Cells.Item(363, 12).Value = "Especial"
$ws.Cells.Item(363, 13).Value = 200
$ws.Cells.Item(363, 14).Value = 20000
$ws.Cells.Item(363, 15).Value = 20000
$ws.Cells.Item(363, 16).Value = 20000
$ws.Cells.Item(363, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(363, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(363, 19).Value = 1333
$ws.Cells.Item(363, 20).Value = 15

# --- Row 364: Nectarín / Super Queen / Primera ---
$ws.Cells.Item(364, 1).Value  = 4
$ws.Cells.Item(364, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(364, 3).Value  = "Los Lagos"
$ws.Cells.Item(364, 4).Value  = 44585
$ws.Cells.Item(364, 5).Value  = 10
$ws.Cells.Item(364, 6).Value  = "Fruta"
$ws.Cells.Item(364, 7).Value  = 100103
$ws.Cells.Item(364, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(364, 9).Value  = 100103006
$ws.Cells.Item(364, 10).Value = "Nectarín"
$ws.Cells.Item(364, 11).Value = "Super Queen"
$ws.Cells.Item(364, 12).Value = "Primera"
$ws.Cells.Item(364, 13).Value = 400
$ws.Cells.Item(364, 14).Value = 16000
$ws.Cells.Item(364, 15).Value = 16000
$ws.Cells.Item(364, 16).Value = 16000
$ws.Cells.Item(364, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(364, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(364, 19).Value = 1067
$ws.Cells.Item(364, 20).Value = 15
